$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-15 03:02:30"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
